$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.725.31"
$ws.Range("E2").Value = "  +3.66%  "
$ws.Range("D3").Value = "1.867.29"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").Value = "231.70"
$ws.Range("E5").Value = "  +2.71%  "
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  +3.76%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "42.84"
$ws.Range("E8").Value = "  +11.49%  "
$ws.Range("E9").Value = "  +7.58%  "
$ws.Range("D10").Value = "0.0698"
$ws.Range("E10").Value = "  +3.55%  "
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("D12").Value = "2.137.47"
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("E13").Value = "  +4.32%  "
$ws.Range("D14").Value = "1.869.21"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("E15").Value = "  +8.23%  "
$ws.Range("D16").Value = "4.76"
$ws.Range("E16").Value = "  +7.62%  "
$ws.Range("D17").Value = "35.751.43"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "70.70"
$ws.Range("E18").Value = "  +3.52%  "
$ws.Range("D19").Value = "249.54"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").Value = "0.0₃0810"
$ws.Range("E20").Value = "  +4.85%  "
$ws.Range("D21").Value = "12.39"
$ws.Range("E21").Value = "  +10.48%  "
$ws.Range("D22").Value = "4.76"
$ws.Range("E22").Value = "  +15.62%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "171.22"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "8.05"
$ws.Range("E26").Value = "  +3.57%  "
$ws.Range("D27").Value = "17.94"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("E28").Value = "  +2.46%  "
$ws.Range("E29").Value = "  +16.36%  "
$ws.Range("D30").Value = "1.01"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").Value = "3.335.68"
$ws.Range("E31").Value = "  +37.29%  "
$ws.Range("D32").Value = "0.0551"
$ws.Range("E32").Value = "  +6.58%  "
$ws.Range("E33").Value = "  +4.89%  "
$ws.Range("D34").Value = "4.09"
$ws.Range("E34").Value = "  +6.31%  "
$ws.Range("E35").Value = "  +5.12%  "
$ws.Range("D36").Value = "99.57"
$ws.Range("E36").Value = "  +21.83%  "
$ws.Range("D37").Value = "0.691"
$ws.Range("E37").Value = "  +7.58%  "
$ws.Range("E38").Value = "  +7.31%  "
$ws.Range("D39").Value = "1.369.91"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("D41").Value = "0.0197"
$ws.Range("E41").Value = "  +6.00%  "
$ws.Range("D42").Value = "1.02"
$ws.Range("E42").Value = "  +7.39%  "
$ws.Range("D43").Value = "14.97"
$ws.Range("E43").Value = "  +8.49%  "
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").Value = "2.84"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("E47").Value = "  +8.96%  "
$ws.Range("D48").Value = "0.0519"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("D49").Value = "2.036.10"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("D50").Value = "105.10"
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("E51").Value = "  +0.37%  "
